$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.7750008106231689
$ws.Cells.Item(2, 5).Value = 219.1298500056891
$ws.Cells.Item(2, 6).Value = 0.007179892215417086
$ws.Cells.Item(2, 7).Value = 0.00607442287755429
$ws.Cells.Item(2, 8).Value = 0.005483430338235971
$ws.Cells.Item(2, 9).Value = 0.005483430338235971
$ws.Cells.Item(2, 10).Value = 0.005170003376302528
$ws.Cells.Item(2, 11).Value = 0.005170003376302528
$ws.Cells.Item(2, 12).Value = 0.005144581724378415
$ws.Cells.Item(2, 13).Value = 0.004898340816651904
$ws.Cells.Item(2, 14).Value = 0.004840728140012362
$ws.Cells.Item(2, 15).Value = 0.004713596928340103
$ws.Cells.Item(2, 16).Value = 0.004629592096009916
$ws.Cells.Item(2, 17).Value = 0.004518444276660413
$ws.Cells.Item(2, 18).Value = 0.004518444276660413
$ws.Cells.Item(2, 19).Value = 0.004417410087141555
$ws.Cells.Item(2, 20).Value = 0.004417410087141555
$ws.Cells.Item(2, 21).Value = 0.00438545263507949
$ws.Cells.Item(2, 22).Value = 0.004311403711929913
$ws.Cells.Item(2, 23).Value = 0.004288767552464166
$ws.Cells.Item(2, 24).Value = 0.004282140208583759
$ws.Cells.Item(2, 25).Value = 0.004271537037147935

$ws.Cells.Item(3, 3).Value = 0.7329962253570557
$ws.Cells.Item(3, 5).Value = 219.602085135708
$ws.Cells.Item(3, 6).Value = 0.006945981785452513
$ws.Cells.Item(3, 7).Value = 0.005476542462444613
$ws.Cells.Item(3, 8).Value = 0.005476542462444613
$ws.Cells.Item(3, 9).Value = 0.005476542462444613
$ws.Cells.Item(3, 10).Value = 0.005246170697833417
$ws.Cells.Item(3, 11).Value = 0.005052969129294182
$ws.Cells.Item(3, 12).Value = 0.004983237214659502
$ws.Cells.Item(3, 13).Value = 0.004830756944608691
$ws.Cells.Item(3, 14).Value = 0.004708135802510487
$ws.Cells.Item(3, 15).Value = 0.004708135802510487
$ws.Cells.Item(3, 16).Value = 0.004645895182482837
$ws.Cells.Item(3, 17).Value = 0.004622096660475992
$ws.Cells.Item(3, 18).Value = 0.004398949352436854
$ws.Cells.Item(3, 19).Value = 0.004398949352436854
$ws.Cells.Item(3, 20).Value = 0.004324790193083907
$ws.Cells.Item(3, 21).Value = 0.004324790193083907
$ws.Cells.Item(3, 22).Value = 0.004322130424624327
$ws.Cells.Item(3, 23).Value = 0.004308949643543376
$ws.Cells.Item(3, 24).Value = 0.004286983255948413
$ws.Cells.Item(3, 25).Value = 0.004280742400306197

$ws.Cells.Item(4, 3).Value = 0.7239871025085449
$ws.Cells.Item(4, 5).Value = 222.5693178281599
$ws.Cells.Item(4, 6).Value = 0.007140213844881211
$ws.Cells.Item(4, 7).Value = 0.006069224252995707
$ws.Cells.Item(4, 8).Value = 0.005610255044257587
$ws.Cells.Item(4, 9).Value = 0.005214629930701066
$ws.Cells.Item(4, 10).Value = 0.005214629930701066
$ws.Cells.Item(4, 11).Value = 0.005042037098638212
$ws.Cells.Item(4, 12).Value = 0.004822007471352777
$ws.Cells.Item(4, 13).Value = 0.004798560290753556
$ws.Cells.Item(4, 14).Value = 0.004793752097648003
$ws.Cells.Item(4, 15).Value = 0.00476679833552344
$ws.Cells.Item(4, 16).Value = 0.004560687105341132
$ws.Cells.Item(4, 17).Value = 0.004560687105341132
$ws.Cells.Item(4, 18).Value = 0.004556136341989701
$ws.Cells.Item(4, 19).Value = 0.004555006623177957
$ws.Cells.Item(4, 20).Value = 0.004470852472555141
$ws.Cells.Item(4, 21).Value = 0.004453860239652294
$ws.Cells.Item(4, 22).Value = 0.004453860239652294
$ws.Cells.Item(4, 23).Value = 0.004363796517886724
$ws.Cells.Item(4, 24).Value = 0.004363796517886724
$ws.Cells.Item(4, 25).Value = 0.004338583193531381

$ws.Cells.Item(5, 3).Value = 0.7590217590332031
$ws.Cells.Item(5, 5).Value = 224.7271922782857
$ws.Cells.Item(5, 6).Value = 0.007179892215417086
$ws.Cells.Item(5, 7).Value = 0.005948076563753606
$ws.Cells.Item(5, 8).Value = 0.005260627160971575
$ws.Cells.Item(5, 9).Value = 0.005260627160971575
$ws.Cells.Item(5, 10).Value = 0.005260627160971575
$ws.Cells.Item(5, 11).Value = 0.005227821085808026
$ws.Cells.Item(5, 12).Value = 0.004920649726104774
$ws.Cells.Item(5, 13).Value = 0.004828518542607388
$ws.Cells.Item(5, 14).Value = 0.004769368717652309
$ws.Cells.Item(5, 15).Value = 0.004754705565014516
$ws.Cells.Item(5, 16).Value = 0.004640500072181872
$ws.Cells.Item(5, 17).Value = 0.00459385493552119
$ws.Cells.Item(5, 18).Value = 0.004554862218794075
$ws.Cells.Item(5, 19).Value = 0.004554862218794075
$ws.Cells.Item(5, 20).Value = 0.004517338553928728
$ws.Cells.Item(5, 21).Value = 0.004510804293041485
$ws.Cells.Item(5, 22).Value = 0.004486580950692511
$ws.Cells.Item(5, 23).Value = 0.004449582558139468
$ws.Cells.Item(5, 24).Value = 0.00438064702296853
$ws.Cells.Item(5, 25).Value = 0.00438064702296853

$ws.Cells.Item(6, 3).Value = 0.8240118026733398
$ws.Cells.Item(6, 5).Value = 222.0081262425665
$ws.Cells.Item(6, 6).Value = 0.006905408613762955
$ws.Cells.Item(6, 7).Value = 0.005940872451807314
$ws.Cells.Item(6, 8).Value = 0.005286937599900098
$ws.Cells.Item(6, 9).Value = 0.005053092989104977
$ws.Cells.Item(6, 10).Value = 0.005046603609362579
$ws.Cells.Item(6, 11).Value = 0.00486501256702489
$ws.Cells.Item(6, 12).Value = 0.00486501256702489
$ws.Cells.Item(6, 13).Value = 0.004665114416238332
$ws.Cells.Item(6, 14).Value = 0.004620661858181286
$ws.Cells.Item(6, 15).Value = 0.004592586815958984
$ws.Cells.Item(6, 16).Value = 0.004592586815958984
$ws.Cells.Item(6, 17).Value = 0.004592586815958984
$ws.Cells.Item(6, 18).Value = 0.004543469278613658
$ws.Cells.Item(6, 19).Value = 0.004519737547145789
$ws.Cells.Item(6, 20).Value = 0.004486579576314326
$ws.Cells.Item(6, 21).Value = 0.004449383001022487
$ws.Cells.Item(6, 22).Value = 0.004418218904881575
$ws.Cells.Item(6, 23).Value = 0.004381863100337099
$ws.Cells.Item(6, 24).Value = 0.004366999538190064
$ws.Cells.Item(6, 25).Value = 0.004327643786404804

$ws.Cells.Item(7, 3).Value = 0.7579965591430664
$ws.Cells.Item(7, 5).Value = 221.8984880090284
$ws.Cells.Item(7, 6).Value = 0.006849836394723634
$ws.Cells.Item(7, 7).Value = 0.005977023322952144
$ws.Cells.Item(7, 8).Value = 0.005641496864181565
$ws.Cells.Item(7, 9).Value = 0.005561738236813441
$ws.Cells.Item(7, 10).Value = 0.005434880373669553
$ws.Cells.Item(7, 11).Value = 0.004876236048119636
$ws.Cells.Item(7, 12).Value = 0.004876236048119636
$ws.Cells.Item(7, 13).Value = 0.004846479835872436
$ws.Cells.Item(7, 14).Value = 0.004713713007143891
$ws.Cells.Item(7, 15).Value = 0.004559276829340798
$ws.Cells.Item(7, 16).Value = 0.004549390763049724
$ws.Cells.Item(7, 17).Value = 0.004549390763049724
$ws.Cells.Item(7, 18).Value = 0.004543994489112092
$ws.Cells.Item(7, 19).Value = 0.004422528246314649
$ws.Cells.Item(7, 20).Value = 0.004422528246314649
$ws.Cells.Item(7, 21).Value = 0.004422528246314649
$ws.Cells.Item(7, 22).Value = 0.0044147526784698
$ws.Cells.Item(7, 23).Value = 0.004395727037281105
$ws.Cells.Item(7, 24).Value = 0.004347938849371638
$ws.Cells.Item(7, 25).Value = 0.004325506588869948

$ws.Cells.Item(8, 3).Value = 0.9825427532196045
$ws.Cells.Item(8, 5).Value = 226.0407983556415
$ws.Cells.Item(8, 6).Value = 0.006793011873349341
$ws.Cells.Item(8, 7).Value = 0.005991284036954737
$ws.Cells.Item(8, 8).Value = 0.005596642936174102
$ws.Cells.Item(8, 9).Value = 0.005320853738428365
$ws.Cells.Item(8, 10).Value = 0.005301193635980267
$ws.Cells.Item(8, 11).Value = 0.004957853786650304
$ws.Cells.Item(8, 12).Value = 0.004896260009063233
$ws.Cells.Item(8, 13).Value = 0.004896260009063233
$ws.Cells.Item(8, 14).Value = 0.004844525199089713
$ws.Cells.Item(8, 15).Value = 0.004836110595201051
$ws.Cells.Item(8, 16).Value = 0.0047530172124069
$ws.Cells.Item(8, 17).Value = 0.0047530172124069
$ws.Cells.Item(8, 18).Value = 0.004641654247031561
$ws.Cells.Item(8, 19).Value = 0.004524275613786623
$ws.Cells.Item(8, 20).Value = 0.004480591593896532
$ws.Cells.Item(8, 21).Value = 0.004480591593896532
$ws.Cells.Item(8, 22).Value = 0.004480591593896532
$ws.Cells.Item(8, 23).Value = 0.00440625337925227
$ws.Cells.Item(8, 24).Value = 0.00440625337925227
$ws.Cells.Item(8, 25).Value = 0.00440625337925227

$ws.Cells.Item(9, 3).Value = 0.8300018310546875
$ws.Cells.Item(9, 5).Value = 228.0460978830051
$ws.Cells.Item(9, 6).Value = 0.007100698031970441
$ws.Cells.Item(9, 7).Value = 0.005998565036348722
$ws.Cells.Item(9, 8).Value = 0.005794980924414307
$ws.Cells.Item(9, 9).Value = 0.005186297634275523
$ws.Cells.Item(9, 10).Value = 0.005186297634275523
$ws.Cells.Item(9, 11).Value = 0.005140199196278499
$ws.Cells.Item(9, 12).Value = 0.004964384731699174
$ws.Cells.Item(9, 13).Value = 0.004964384731699174
$ws.Cells.Item(9, 14).Value = 0.004834289355263715
$ws.Cells.Item(9, 15).Value = 0.004834289355263715
$ws.Cells.Item(9, 16).Value = 0.004691388340106782
$ws.Cells.Item(9, 17).Value = 0.004608427495743765
$ws.Cells.Item(9, 18).Value = 0.004587609157411227
$ws.Cells.Item(9, 19).Value = 0.004587609157411227
$ws.Cells.Item(9, 20).Value = 0.004587609157411227
$ws.Cells.Item(9, 21).Value = 0.004514625665163507
$ws.Cells.Item(9, 22).Value = 0.004514625665163507
$ws.Cells.Item(9, 23).Value = 0.004514625665163507
$ws.Cells.Item(9, 24).Value = 0.00449348331038414
$ws.Cells.Item(9, 25).Value = 0.00444534303865507

$ws.Cells.Item(10, 3).Value = 0.7560064792633057
$ws.Cells.Item(10, 5).Value = 219.085237001429
$ws.Cells.Item(10, 6).Value = 0.007049674773747626
$ws.Cells.Item(10, 7).Value = 0.006230786531379801
$ws.Cells.Item(10, 8).Value = 0.005500396363843224
$ws.Cells.Item(10, 9).Value = 0.005422556020622016
$ws.Cells.Item(10, 10).Value = 0.005286029712903918
$ws.Cells.Item(10, 11).Value = 0.004979645791404658
$ws.Cells.Item(10, 12).Value = 0.004979645791404658
$ws.Cells.Item(10, 13).Value = 0.004915234714416546
$ws.Cells.Item(10, 14).Value = 0.004770396603170345
$ws.Cells.Item(10, 15).Value = 0.004653188370405039
$ws.Cells.Item(10, 16).Value = 0.004617788898394118
$ws.Cells.Item(10, 17).Value = 0.004455308931067841
$ws.Cells.Item(10, 18).Value = 0.004455308931067841
$ws.Cells.Item(10, 19).Value = 0.004455308931067841
$ws.Cells.Item(10, 20).Value = 0.004455308931067841
$ws.Cells.Item(10, 21).Value = 0.004342063852462543
$ws.Cells.Item(10, 22).Value = 0.004342063852462543
$ws.Cells.Item(10, 23).Value = 0.004286221662520001
$ws.Cells.Item(10, 24).Value = 0.004286221662520001
$ws.Cells.Item(10, 25).Value = 0.004270667387942086

$ws.Cells.Item(11, 3).Value = 0.7239811420440674
$ws.Cells.Item(11, 5).Value = 222.277555050392
$ws.Cells.Item(11, 6).Value = 0.0070113205034393
$ws.Cells.Item(11, 7).Value = 0.005983212000766911
$ws.Cells.Item(11, 8).Value = 0.005759213262856815
$ws.Cells.Item(11, 9).Value = 0.005296487017556103
$ws.Cells.Item(11, 10).Value = 0.005296487017556103
$ws.Cells.Item(11, 11).Value = 0.005151643623829331
$ws.Cells.Item(11, 12).Value = 0.005125584347232126
$ws.Cells.Item(11, 13).Value = 0.004808521021150545
$ws.Cells.Item(11, 14).Value = 0.004635722078536997
$ws.Cells.Item(11, 15).Value = 0.004635722078536997
$ws.Cells.Item(11, 16).Value = 0.004567693541758223
$ws.Cells.Item(11, 17).Value = 0.004470808552454414
$ws.Cells.Item(11, 18).Value = 0.004470808552454414
$ws.Cells.Item(11, 19).Value = 0.004470808552454414
$ws.Cells.Item(11, 20).Value = 0.004464119714581672
$ws.Cells.Item(11, 21).Value = 0.004406925566519193
$ws.Cells.Item(11, 22).Value = 0.004362755924277919
$ws.Cells.Item(11, 23).Value = 0.004362755924277919
$ws.Cells.Item(11, 24).Value = 0.004352575258936626
$ws.Cells.Item(11, 25).Value = 0.004332895809949161
